$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sprint 6 backlog" section heading was actually a duplicate/mislabeled
# "Sprint 3 backlog" -> fix the label
$ws.Range("B62").Value = "Sprint 3 backlog"

# The "Test Webcam" task (row 65) was dropped from the backlog; clear the
# whole row and give it a plain, centered, borderless look
$rng = $ws.Range("B65:H65")
$rng.Value = $null
$rng.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$rng.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# Restore the view/selection that was active when the author saved
[void]$excel.Goto($ws.Range("A43"), $true)
$ws.Range("F64").Select() | Out-Null
